$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-11 17:48:25"
$ws.Range("I2").Value = "2.0 mm"
$ws.Range("O2").Value = "3.2 °C"
$ws.Range("E3").Value = "2026-02-11 17:48:27"
$ws.Range("O3").Value = "0.1 °C"
$ws.Range("E4").Value = "2026-02-11 17:48:30"
$ws.Range("J4").Value = "1002.8 hPa"
$ws.Range("O4").Value = "15.9 °C"
$ws.Range("E5").Value = "2026-02-11 17:48:32"
$ws.Range("I5").Value = "0.8 mm"
$ws.Range("L5").Value = "49.3 km/h - 273º 17:19 TU"
$ws.Range("M5").Value = "4.4 °C 17:20 TU"
$ws.Range("O5").Value = "0.4 °C"
$ws.Range("E6").Value = "2026-02-11 17:48:34"
$ws.Range("J6").Value = "1003.3 hPa"
$ws.Range("E7").Value = "2026-02-11 17:48:37"
$ws.Range("E8").Value = "2026-02-11 17:48:39"
$ws.Range("J8").Value = "1003.4 hPa"
$ws.Range("O8").Value = "15.5 °C"
$ws.Range("E9").Value = "2026-02-11 17:48:42"
$ws.Range("E10").Value = "2026-02-11 17:48:44"
$ws.Range("E11").Value = "2026-02-11 17:48:47"
$ws.Range("H11").Value = "'79%"
$ws.Range("O11").Value = "7.9 °C"
$ws.Range("E12").Value = "2026-02-11 17:48:49"
$ws.Range("O12").Value = "11.6 °C"
$ws.Range("E13").Value = "2026-02-11 17:48:51"
$ws.Range("J13").Value = "1005.4 hPa"
$ws.Range("E14").Value = "2026-02-11 17:48:54"
$ws.Range("E15").Value = "2026-02-11 17:48:56"
$ws.Range("E16").Value = "2026-02-11 17:48:59"
$ws.Range("I16").Value = "4.8 mm"
$ws.Range("E17").Value = "2026-02-11 17:49:01"
$ws.Range("H17").Value = "'75%"
$ws.Range("K17").Value = "10.7 MJ/m2"
$ws.Range("E18").Value = "2026-02-11 17:49:04"
$ws.Range("H18").Value = "'72%"
$ws.Range("J18").Value = "1003.3 hPa"
$ws.Range("L18").Value = "24.5 km/h - 263º 17:19 TU"
$ws.Range("O18").Value = "13.9 °C"
$ws.Range("E19").Value = "2026-02-11 17:49:06"
$ws.Range("H19").Value = "'78%"
$ws.Range("E20").Value = "2026-02-11 17:49:08"
$ws.Range("I20").Value = "0.4 mm"
$ws.Range("E21").Value = "2026-02-11 17:49:11"
$ws.Range("I21").Value = "1.4 mm"
$ws.Range("J21").Value = "1005.9 hPa"
$ws.Range("O21").Value = "8.5 °C"
$ws.Range("E22").Value = "2026-02-11 17:49:13"
$ws.Range("G22").Value = "120 cm"
$ws.Range("M22").Value = "-0.7 °C 17:16 TU"
$ws.Range("E23").Value = "2026-02-11 17:49:16"
$ws.Range("H23").Value = "'70%"
$ws.Range("I23").Value = "3.2 mm"
$ws.Range("E24").Value = "2026-02-11 17:49:18"
$ws.Range("H24").Value = "'72%"
$ws.Range("I24").Value = "5.9 mm"
$ws.Range("J24").Value = "1007.3 hPa"
$ws.Range("N24").Value = "11.2 °C 17:01 TU"
$ws.Range("O24").Value = "13.5 °C"
$ws.Range("E25").Value = "2026-02-11 17:49:21"
$ws.Range("L25").Value = "51.5 km/h - 267º 17:23 TU"
$ws.Range("E26").Value = "2026-02-11 17:49:23"
$ws.Range("E27").Value = "2026-02-11 17:49:26"
$ws.Range("I27").Value = "0.9 mm"
$ws.Range("O27").Value = "0.5 °C"
$ws.Range("E28").Value = "2026-02-11 17:49:28"
$ws.Range("J28").Value = "1003.6 hPa"
$ws.Range("O28").Value = "10.8 °C"
$ws.Range("E29").Value = "2026-02-11 17:49:31"
$ws.Range("E30").Value = "2026-02-11 17:49:33"
$ws.Range("J30").Value = "1003.5 hPa"
$ws.Range("E31").Value = "2026-02-11 17:49:36"
$ws.Range("J31").Value = "1002.7 hPa"
$ws.Range("E32").Value = "2026-02-11 17:49:38"
$ws.Range("E33").Value = "2026-02-11 17:49:40"
$ws.Range("J33").Value = "1005.0 hPa"
$ws.Range("O33").Value = "6.8 °C"
$ws.Range("E34").Value = "2026-02-11 17:49:43"
$ws.Range("O34").Value = "3.5 °C"
$ws.Range("E35").Value = "2026-02-11 17:49:45"
$ws.Range("J35").Value = "1008.1 hPa"
$ws.Range("K35").Value = "7.3 MJ/m2"
$ws.Range("E36").Value = "2026-02-11 17:49:48"
$ws.Range("H36").Value = "'87%"
$ws.Range("J36").Value = "1003.6 hPa"
$ws.Range("E37").Value = "2026-02-11 17:49:50"
$ws.Range("H37").Value = "'82%"
$ws.Range("J37").Value = "1004.8 hPa"
$ws.Range("E38").Value = "2026-02-11 17:49:53"
$ws.Range("E39").Value = "2026-02-11 17:49:55"
$ws.Range("L39").Value = "76.3 km/h - 301º 17:03 TU"
$ws.Range("E40").Value = "2026-02-11 17:49:58"
$ws.Range("I40").Value = "1.3 mm"
$ws.Range("J40").Value = "1007.3 hPa"
$ws.Range("E41").Value = "2026-02-11 17:50:00"
$ws.Range("H41").Value = "'44%"
$ws.Range("J41").Value = "1005.1 hPa"
$ws.Range("N41").Value = "16.9 °C 17:03 TU"
$ws.Range("O41").Value = "19.3 °C"
$ws.Range("E42").Value = "2026-02-11 17:50:05"
$ws.Range("E43").Value = "2026-02-11 17:50:07"
$ws.Range("E44").Value = "2026-02-11 17:50:10"
$ws.Range("I44").Value = "4.6 mm"
$ws.Range("E45").Value = "2026-02-11 17:50:12"
$ws.Range("I45").Value = "1.1 mm"
$ws.Range("J45").Value = "1006.2 hPa"
$ws.Range("O45").Value = "6.8 °C"
$ws.Range("E46").Value = "2026-02-11 17:50:15"
$ws.Range("H46").Value = "'55%"
$ws.Range("J46").Value = "1007.8 hPa"
$ws.Range("N46").Value = "13.1 °C 17:29 TU"
$ws.Range("O46").Value = "17.5 °C"

Write-Host "Applied all updates"
